$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 10 (pushes the existing rows 10-26 down to 11-27).
$ws.Rows(10).Insert()
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44811
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112037
$ws.Range("G10").Value = "Cebollín"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("N10").Value = "$/docena de atados"
$ws.Range("O10").Value = "Provincia de Diguillín"
$ws.Range("P10").Value = 2833
$ws.Range("Q10").Value = 3
$ws.Range("R10").Value = "Hortaliza"

# Insert a second new data row at row 18 (after the first insert, the row that
# held 2022-10-22 / 44790 is now row 19; pushes rows 18-27 down to 19-28).
$ws.Rows(18).Insert()
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44810
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112037
$ws.Range("G18").Value = "Cebollín"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8500
$ws.Range("N18").Value = "$/docena de atados"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 2833
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = "Hortaliza"
